$d = $word.ActiveDocument

# b. ... so we need to create a Q-Q plot to be sure  ->  ... a histogram to be sure
$d.Content.Find.Execute(
    "Q-Q plot to be sure", $false, $false, $false, $false, $false,
    $true, 1, $false, "histogram to be sure", 2)

# the distribution of sample means is normal:
#   -> the distribution of differences is normal in order to be able to claim
#      that the distribution sample means is normal:
$d.Content.Find.Execute(
    "the distribution of sample means is normal:", $false, $false, $false, $false, $false,
    $true, 1, $false,
    "the distribution of differences is normal in order to be able to claim that the distribution sample means is normal:",
    2)

# The data appear to be normal  ->  The differences appear to be normal
$d.Content.Find.Execute(
    "The data appear to be normal", $false, $false, $false, $false, $false,
    $true, 1, $false, "The differences appear to be normal", 2)
